$wb = $excel.ActiveWorkbook

$lora   = $wb.Worksheets.Item("LoRa")
$sigfox = $wb.Worksheets.Item("Sigfox")

# --- LoRa sheet: add the "Hexa 8 octects" column to the small identifiers table ---

# Header cell, formatted like the other header cells on that row (A8/B8).
$lora.Range("C8").Value = "Hexa 8 octects"
$lora.Range("A8").Copy()
$lora.Range("C8").PasteSpecial(-4122)  # xlPasteFormats

# Data cells - set the values first, then copy the row's look (fill/border)
# from column A, and finally force the "Text" number format used for the
# other Hexa-8-octets values in this workbook (see C3/C4).
$lora.Range("C9").Value  = "50f925fb"
$lora.Range("C10").Value = 64754638

$lora.Range("A9:A10").Copy()
$lora.Range("C9:C10").PasteSpecial(-4122)  # xlPasteFormats
$lora.Range("C9:C10").NumberFormat = "@"

# --- Selection / active-sheet bookkeeping ---
# Sigfox is no longer the active tab; its selection moves to G12.
$sigfox.Range("G12").Select()

# LoRa becomes the active sheet, with E10 selected.
$lora.Activate()
$lora.Range("E10").Select()
